$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row directly below the current row 8 ("Upload"). This
# shifts row 9 ("Desc" header) and everything after it down by one, while
# leaving row 8 itself (and its row-level formatting) untouched.
$ws.Rows.Item(9).Insert()

# Move the "Upload" row's values down into the newly blank row 9, and copy
# row 8's formatting down onto it too (it's currently blank/default).
$ws.Range("A9").Value2 = $ws.Range("A8").Value2
$ws.Range("B9").Value2 = $ws.Range("B8").Value2
$ws.Range("C9").Value2 = $ws.Range("C8").Value2
$ws.Range("D9").Value2 = $ws.Range("D8").Value2
$ws.Range("E9").Value2 = $ws.Range("E8").Value2

$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)

# Now overwrite row 8 with the new "Force" row contents.
$ws.Range("A8").Value = "Force"
$ws.Range("B8").Value = $false
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $false
